$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Delete the duplicate "Contact / No display for ContactDetail" row (row 11)
# so the remaining Contact row (row 10) can be repurposed.
$ws.Rows(11).Delete()

# Update version and date metadata values.
$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty; now populated.
$ws.Range("B9").Value = "Alvearie Team"

# Former "Contact" / "No display for ContactDetail" row becomes a
# "Jurisdiction" / "United States of America" row.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
